$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 320
$ws.Range("F3").Value = 1111
$ws.Range("F5").Value = 1118
$ws.Range("F6").Value = 3370
$ws.Range("F7").Value = 60
$ws.Range("F9").Value = 1179
$ws.Range("F10").Value = 760
$ws.Range("F11").Value = 587
$ws.Range("F13").Value = 53
$ws.Range("F14").Value = 146
$ws.Range("F16").Value = 1755
$ws.Range("F17").Value = 40
$ws.Range("F18").Value = 353
$ws.Range("F19").Value = 30
$ws.Range("F20").Value = 49
$ws.Range("F21").Value = 652
$ws.Range("F22").Value = 398
$ws.Range("F23").Value = 714
$ws.Range("F24").Value = 78882
$ws.Range("F25").Value = 78882
$ws.Range("F26").Value = 704
$ws.Range("F27").Value = 655
$ws.Range("F28").Value = 33596
$ws.Range("F29").Value = 33596
$ws.Range("F30").Value = 496
$ws.Range("F32").Value = 18
$ws.Range("F34").Value = 31
$ws.Range("F35").Value = 957
$ws.Range("F36").Value = 284
$ws.Range("F37").Value = 159
$ws.Range("F38").Value = 570
$ws.Range("F39").Value = 899
$ws.Range("F40").Value = 1181
$ws.Range("F41").Value = 5451
$ws.Range("F42").Value = 761
$ws.Range("F44").Value = 1
$ws.Range("F46").Value = 378
$ws.Range("F50").Value = 44

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 23
$ws.Range("F8").Value = 18
$ws.Range("F15").Value = 1710
$ws.Range("F16").Value = 9
$ws.Range("F21").Value = 66
$ws.Range("F24").Value = 512
$ws.Range("F25").Value = 512
$ws.Range("F26").Value = 11
$ws.Range("F27").Value = 765
$ws.Range("F28").Value = 13
$ws.Range("F35").Value = 1662
$ws.Range("F36").Value = 492
$ws.Range("F47").Value = 63
$ws.Range("F48").Value = 823
$ws.Range("F49").Value = 135
$ws.Range("F50").Value = 45

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F4").Value = 725
$ws.Range("F5").Value = 558
$ws.Range("F6").Value = 589
$ws.Range("F7").Value = 78

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 725
$ws.Range("F3").Value = 320
$ws.Range("F4").Value = 558
$ws.Range("F5").Value = 1111
$ws.Range("F7").Value = 1118
$ws.Range("F8").Value = 3370
$ws.Range("F9").Value = 60
$ws.Range("F10").Value = 1179
$ws.Range("F11").Value = 760
$ws.Range("F12").Value = 589
$ws.Range("F13").Value = 589
$ws.Range("F16").Value = 53
$ws.Range("F17").Value = 146
$ws.Range("F19").Value = 1755
$ws.Range("F20").Value = 353
$ws.Range("F22").Value = 30
$ws.Range("F23").Value = 49
$ws.Range("F24").Value = 652
$ws.Range("F26").Value = 398
$ws.Range("F27").Value = 714
$ws.Range("F28").Value = 78882
$ws.Range("F29").Value = 655
$ws.Range("F30").Value = 33596
$ws.Range("F31").Value = 496
$ws.Range("F33").Value = 512
$ws.Range("F34").Value = 31
$ws.Range("F35").Value = 11
$ws.Range("F36").Value = 957
$ws.Range("F37").Value = 13
$ws.Range("F38").Value = 284
$ws.Range("F39").Value = 570
$ws.Range("F40").Value = 899
$ws.Range("F41").Value = 899
$ws.Range("F42").Value = 1181
$ws.Range("F43").Value = 5451
$ws.Range("F44").Value = 761
$ws.Range("F45").Value = 1662
$ws.Range("F46").Value = 492
$ws.Range("F49").Value = 378
$ws.Range("F51").Value = 63
$ws.Range("F53").Value = 823
$ws.Range("F54").Value = 136
$ws.Range("F55").Value = 44
